$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 19

$ws.Cells.Item($row, 1).Value = 42601.898923611108
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($row, 2).Value = "Noun"

$ws.Cells.Item($row, 3).Value = 13599
$ws.Cells.Item($row, 4).Value = 8966
$ws.Cells.Item($row, 5).Value = 1638
$ws.Cells.Item($row, 6).Value = 179
$ws.Cells.Item($row, 7).Value = 81
$ws.Cells.Item($row, 8).Value = 68
$ws.Cells.Item($row, 9).Value = 30
$ws.Cells.Item($row, 10).Value = 2
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 100
$ws.Cells.Item($row, 13).Value = 0
